$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2/C2 (Elkind -> Gribnich)
$ws.Range("B2").Value = "Gribnich"
$ws.Range("C2").Value = "sgribnich345@gmail.com"

# Update B3/C3 (Elkind -> Smith)
$ws.Range("B3").Value = "Smith"
$ws.Range("C3").Value = "ssmitty_do_not_use.@aol.com"

# Update D4 topics and C4 email
$ws.Range("D4").Value = "Ukraine, Micropython, Tesla, Subaru, Lindsey Stirling, Raspberry Pi"
$ws.Range("C4").Value = "jlucas.foobar@yahoo.com"

# Column width adjustments (values chosen so the engine's pixel-quantized
# ColumnWidth -> stored "width" matches the target OOXML width exactly)
$ws.Range("C1").ColumnWidth = 27.5
$ws.Range("D1").ColumnWidth = 60.15

# Selection change
$ws.Range("B9").Select() | Out-Null
